$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$ws.Range("A17").Value = "WAT"
$ws.Range("B17").Value = "WoS Author Transformation testcases"
$ws.Range("C17").Value = "Y"

# Copy the style/formatting from the row above (A16:C16) down to the new row (A17:C17)
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F13").Select() | Out-Null
